$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21
$ws.Range("A21").Value = 42610
$ws.Range("A21").NumberFormat = "d-mmm"
$ws.Range("B21").Value = "think up new UI to enhance chords"
$ws.Range("B21").Interior.Color = 65535
$ws.Range("C21").Value = "Pending"

# Row 22
$ws.Range("A22").Value = 42610
$ws.Range("A22").NumberFormat = "d-mmm"
$ws.Range("B22").Value = "adjust sharp/flat based on melody note and key"
$ws.Range("B22").Interior.Color = 65535
$ws.Range("C22").Value = "Resolved"

$ws.Range("B22").Select()
